$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = -0.07875464121263161
$ws.Range("C2").Value = 0.2712360448979064
$ws.Range("D2").Value = 137.6354536320966
$ws.Range("E2").Value = -0.2903546290917011
$ws.Range("F2").Value = 0.7719819679278436
$ws.Range("B3").Value = 0.4729630458858647
$ws.Range("C3").Value = 0.484497005678142
$ws.Range("D3").Value = 166.3444060592113
$ws.Range("E3").Value = 0.976193950309077
$ws.Range("F3").Value = 0.3303858684956607
$ws.Range("B4").Value = -0.3393121207190132
$ws.Range("C4").Value = 0.3960722668307303
$ws.Range("D4").Value = 176.3560324640498
$ws.Range("E4").Value = -0.8566924501786066
$ws.Range("F4").Value = 0.3927774658033908
$ws.Range("B5").Value = 0.1705593857111773
$ws.Range("C5").Value = 0.3851890032934916
$ws.Range("D5").Value = 179.7332719118628
$ws.Range("E5").Value = 0.4427940160618267
$ws.Range("F5").Value = 0.6584471780563337
$ws.Range("B6").Value = 0.6559175507094906
$ws.Range("C6").Value = 0.3258268676812087
$ws.Range("D6").Value = 126.8039750548492
$ws.Range("E6").Value = 2.013086137976949
$ws.Range("F6").Value = 0.0462229974896117
$ws.Range("B7").Value = -0.6231150040980329
$ws.Range("C7").Value = 0.567647069826566
$ws.Range("D7").Value = 273.7431405467593
$ws.Range("E7").Value = -1.097715529983118
$ws.Range("F7").Value = 0.2732933432856151
$ws.Range("B8").Value = -0.5955134525019619
$ws.Range("C8").Value = 0.6478193449045511
$ws.Range("D8").Value = 181.3718460970453
$ws.Range("E8").Value = -0.9192585204285681
$ws.Range("F8").Value = 0.3591812996625936
$ws.Range("B9").Value = -0.4143208248577895
$ws.Range("C9").Value = 0.7009852699593431
$ws.Range("D9").Value = 225.109645999053
$ws.Range("E9").Value = -0.591054965936474
$ws.Range("F9").Value = 0.555076586937634
$ws.Range("B10").Value = -0.9723398683429955
$ws.Range("C10").Value = 0.5507690967395207
$ws.Range("D10").Value = 155.1944505260827
$ws.Range("E10").Value = -1.765421978282945
$ws.Range("F10").Value = 0.07945951067312214
$ws.Range("B11").Value = -0.0625419187438076
$ws.Range("C11").Value = 0.4875380440167622
$ws.Range("D11").Value = 162.8521523765324
$ws.Range("E11").Value = -0.1282811044416819
$ws.Range("F11").Value = 0.8980847975208431
$ws.Range("B12").Value = -0.2715660407377426
$ws.Range("C12").Value = 0.4739038328912701
$ws.Range("D12").Value = 161.3567688706566
$ws.Range("E12").Value = -0.5730403974176111
$ws.Range("F12").Value = 0.5674150503766027
$ws.Range("B13").Value = -0.4311303428920717
$ws.Range("C13").Value = 0.7477125507913269
$ws.Range("D13").Value = 274.4925499811952
$ws.Range("E13").Value = -0.5765990452290701
$ws.Range("F13").Value = 0.5646828655198403
$ws.Range("B14").Value = 0.3709229407557811
$ws.Range("C14").Value = 0.8064446272889852
$ws.Range("D14").Value = 272.6384218948076
$ws.Range("E14").Value = 0.4599484306848296
$ws.Range("F14").Value = 0.6459198154282715
$ws.Range("B15").Value = -0.2226293547218639
$ws.Range("C15").Value = 0.7038126702403158
$ws.Range("D15").Value = 274.9404147957764
$ws.Range("E15").Value = -0.3163190492803254
$ws.Range("F15").Value = 0.7520003699713125
$ws.Range("B16").Value = -0.7326867164159139
$ws.Range("C16").Value = 0.6613422261568976
$ws.Range("D16").Value = 274.9987953100372
$ws.Range("E16").Value = -1.107878322957848
$ws.Range("F16").Value = 0.2688824127433405
$ws.Range("B17").Value = 1.46640135377371
$ws.Range("C17").Value = 0.7552960626962272
$ws.Range("D17").Value = 165.8580198983113
$ws.Range("E17").Value = 1.941492119711318
$ws.Range("F17").Value = 0.05389336822884208
$ws.Range("B18").Value = 0.6567541385582069
$ws.Range("C18").Value = 0.8091813493721264
$ws.Range("D18").Value = 204.3790638140996
$ws.Range("E18").Value = 0.8116278743545023
$ws.Range("F18").Value = 0.4179494309640903
$ws.Range("B19").Value = 0.3778041657856547
$ws.Range("C19").Value = 0.9860342500620759
$ws.Range("D19").Value = 274.9913663094633
$ws.Range("E19").Value = 0.3831552157157523
$ws.Range("F19").Value = 0.7019007036554604
$ws.Range("B20").Value = 0.2917074250094386
$ws.Range("C20").Value = 1.092927224733871
$ws.Range("D20").Value = 270.1808043850413
$ws.Range("E20").Value = 0.26690471095225
$ws.Range("F20").Value = 0.7897461331182813
$ws.Range("B21").Value = 1.06799984748068
$ws.Range("C21").Value = 0.8596073415151869
$ws.Range("D21").Value = 274.9916372831871
$ws.Range("E21").Value = 1.242427554885897
$ws.Range("F21").Value = 0.2151372962932173
$ws.Range("B22").Value = 0.2053956671656148
$ws.Range("C22").Value = 0.9405643817780359
$ws.Range("D22").Value = 274.8690356230804
$ws.Range("E22").Value = 0.2183749152581521
$ws.Range("F22").Value = 0.8272990508515825
$ws.Range("B23").Value = 0.7793469430242418
$ws.Range("C23").Value = 0.8736593374445272
$ws.Range("D23").Value = 274.1285946515516
$ws.Range("E23").Value = 0.8920490053982011
$ws.Range("F23").Value = 0.3731490160142684
$ws.Range("B24").Value = -1.424709943635549
$ws.Range("C24").Value = 1.178607177466117
$ws.Range("D24").Value = 273.8392916307839
$ws.Range("E24").Value = -1.208808134614051
$ws.Range("F24").Value = 0.2277792798609498
$ws.Range("B25").Value = -1.247841531917578
$ws.Range("C25").Value = 1.304629261775828
$ws.Range("D25").Value = 274.621892121622
$ws.Range("E25").Value = -0.9564721323351644
$ws.Range("F25").Value = 0.3396749800492365
